$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.015565
$ws.Range("H2").Value = 0.046695
$ws.Range("I2").Value = 0.4266916434413122
$ws.Range("J2").Value = 0.4266916434413122
$ws.Range("M2").Value = 5.347438999999999
$ws.Range("N2").Value = 16.042317
$ws.Range("O2").Value = 0.1410514228841643
$ws.Range("P2").Value = 0.1410514228841643
$ws.Range("Q2").Value = 0.08323288803499998
$ws.Range("R2").Value = 0.7490959923149999
$ws.Range("S2").Value = 0.06018546344017957
$ws.Range("T2").Value = 0.06018546344017956

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.015565
$ws.Range("H3").Value = 0.046695
$ws.Range("I3").Value = 0.4266916434413122
$ws.Range("J3").Value = 0.4266916434413122
$ws.Range("M3").Value = 24.52471933333333
$ws.Range("N3").Value = 73.574158
$ws.Range("O3").Value = 0.6468978061837527
$ws.Range("P3").Value = 0.6468978061837526
$ws.Range("Q3").Value = 0.3817272564233333
$ws.Range("R3").Value = 3.43554530781
$ws.Range("S3").Value = 0.2760258880591249
$ws.Range("T3").Value = 0.2760258880591248

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.015565
$ws.Range("H4").Value = 0.046695
$ws.Range("I4").Value = 0.4266916434413122
$ws.Range("J4").Value = 0.4266916434413122
$ws.Range("O4").Value = 0.2120507709320831
$ws.Range("P4").Value = 0.2120507709320831
$ws.Range("Q4").Value = 0.1251288197866667
$ws.Range("R4").Value = 1.12615937808
$ws.Range("S4").Value = 0.0904802919420078
$ws.Range("T4").Value = 0.09048029194200777

# Row 5
$ws.Range("I5").Value = 0.5733083565586878
$ws.Range("J5").Value = 0.5733083565586878
$ws.Range("M5").Value = 5.347438999999999
$ws.Range("N5").Value = 16.042317
$ws.Range("O5").Value = 0.1410514228841643
$ws.Range("P5").Value = 0.1410514228841643
$ws.Range("Q5").Value = 0.1118327742866667
$ws.Range("R5").Value = 1.00649496858
$ws.Range("S5").Value = 0.08086595944398471
$ws.Range("T5").Value = 0.08086595944398471

# Row 6
$ws.Range("I6").Value = 0.5733083565586878
$ws.Range("J6").Value = 0.5733083565586878
$ws.Range("M6").Value = 24.52471933333333
$ws.Range("N6").Value = 73.574158
$ws.Range("O6").Value = 0.6468978061837527
$ws.Range("P6").Value = 0.6468978061837526
$ws.Range("Q6").Value = 0.5128936303244445
$ws.Range("R6").Value = 4.61604267292
$ws.Range("S6").Value = 0.3708719181246278
$ws.Range("T6").Value = 0.3708719181246278

# Row 7
$ws.Range("I7").Value = 0.5733083565586878
$ws.Range("J7").Value = 0.5733083565586878
$ws.Range("O7").Value = 0.2120507709320831
$ws.Range("P7").Value = 0.2120507709320831
$ws.Range("S7").Value = 0.1215704789900754
$ws.Range("T7").Value = 0.1215704789900753
